$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Re-style the surviving cells BEFORE the trailing columns are removed,
#     borrowing the exact formatting that is about to be vacated:
#       D5:G5 -> the style already on H5:J5 (thin, no extra border)
#       E6:J6 -> the style already on M6:P6 (medium top+bottom border)
$ws.Range("H5").Copy()
$ws.Range("D5:G5").PasteSpecial(-4122)
$ws.Range("M6").Copy()
$ws.Range("E6:J6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Remove the now-unused trailing columns K:P entirely so the sheet extent
# becomes A1:J6.
$ws.Range("K1:P6").EntireColumn.Delete()

# --- Row 4 (year headers): shift the 13-year run 2007-2019 to the new
#     7-year run 2015-2021 that now lives in D4:J4.
$ws.Range("D4").Value = 2015
$ws.Range("E4").Value = 2016
$ws.Range("F4").Value = 2017
$ws.Range("G4").Value = 2018
$ws.Range("H4").Value = 2019
$ws.Range("I4").Value = 2020
$ws.Range("J4").Value = 2021

# --- Row 5 new values
$ws.Range("D5").Value = 2.2197193775563164
$ws.Range("E5").Value = 2.1235271668715399
$ws.Range("F5").Value = 2.7818537161298167
$ws.Range("G5").Value = 6.7272960584548969
$ws.Range("H5").Value = 5.1525830614767187
$ws.Range("I5").Value = 4.4774536255935971
$ws.Range("J5").Value = 4.6024666695867751

# --- Row 6 new values
$ws.Range("D6").Value = 2.2322863217945752
$ws.Range("E6").Value = 2.8603553109638966
$ws.Range("F6").Value = 3.113207036164539
$ws.Range("G6").Value = 6.2970593463100784
$ws.Range("H6").Value = 4.8617746111834492
$ws.Range("I6").Value = 2.6715092780025032
$ws.Range("J6").Value = 4.3694509108608912

# Column widths for the surviving data columns D:J; everything past column J
# reverts to the default width. (8.6 characters is this engine's closest
# representable input to the authored 9.42578125 stored width.)
$ws.Range("D1:J1").EntireColumn.ColumnWidth = 8.6

# Move/park the active selection the way the saved workbook view shows it.
$ws.Range("K16").Select()
